$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.108.51"
$ws.Range("E2").Value = "  -0.79%  "

# Row 3
$ws.Range("D3").Value = "1.649.84"
$ws.Range("E3").Value = "  -0.92%  "

# Row 4
$ws.Range("E4").Value = "  -0.54%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.80%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5193"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.91%  "

# Row 7
$ws.Range("E7").Value = "  -0.43%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2615"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.70%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06278"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.36%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07784"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.56%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.668.68"
$ws.Range("E12").Value = "  +0.42%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.464"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.23%  "

# Row 14
$ws.Range("D14").Value = "1.875.79"
$ws.Range("E14").Value = "  -0.93%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5542"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.22%  "

# Row 16
$ws.Range("D16").Value = "0.0₅7973"
$ws.Range("E16").Value = "  -3.31%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.70%  "

# Row 18
$ws.Range("D18").Value = "26.087.33"
$ws.Range("E18").Value = "  -0.94%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.40%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.624"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.61%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.10%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.90%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.937"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.74%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.006"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.48%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.28%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1205"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.13%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.157"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.60%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.87%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.476"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.60%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05617"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.91%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.264"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.45%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.473"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.02%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.378"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.95%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.595"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.56%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.800"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.87%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9473"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.26%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.403"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.69%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5641"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.07%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.947"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.34%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01575"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.85%  "

# Row 41
$ws.Range("D41").Value = "1.060.63"
$ws.Range("E41").Value = "  +0.75%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8377"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.76%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.32%  "

# Row 45
$ws.Range("B45").Value = "PAXGold"
$ws.Range("C45").Value = "https://coinranking.com/coin/YRTkUcMi+paxgold-paxg"
$ws.Range("D45").Value = "1.912.41"
$ws.Range("E45").Value = "  -0.38%  "

# Row 46
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.787.15"
$ws.Range("E46").Value = "  -0.95%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.48%  "

# Row 48
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₈106"
$ws.Range("E48").Value = "  -1.12%  "

# Row 49
$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.003"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.04%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05316"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.92%  "

# Row 51
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4334"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.20%  "
